$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1000
$ws.Range("J32").Value = 1000
$ws.Range("L32").Value = 1000
$ws.Range("N32").Value = -1652

$ws.Range("H76").Value = 7216
$ws.Range("I76").Value = 7133
$ws.Range("J76").Value = 7299
$ws.Range("K76").Value = 7133
$ws.Range("L76").Value = 7299
$ws.Range("M76").Value = -6818
$ws.Range("N76").Value = -7929

$ws.Range("H79").Value = 7216
$ws.Range("I79").Value = 7133
$ws.Range("J79").Value = 7299
$ws.Range("K79").Value = 7133
$ws.Range("L79").Value = 7299
$ws.Range("M79").Value = -6041
$ws.Range("N79").Value = -9483

$ws.Range("H112").Value = 2458.5
$ws.Range("I112").Value = 2308.5
$ws.Range("J112").Value = 2488.5
$ws.Range("K112").Value = 6925.5
$ws.Range("L112").Value = 7465.5
$ws.Range("M112").Value = -5817.5
$ws.Range("N112").Value = -9681.5

$ws.Range("H113").Value = 7999.6665
$ws.Range("J113").Value = 7999.6665
$ws.Range("L113").Value = 7999.6665
$ws.Range("N113").Value = -14507.6665

$ws.Range("H127").Value = 118541.53
$ws.Range("I127").Value = 125716.06
$ws.Range("K127").Value = 377148.18
$ws.Range("M127").Value = -372188.18

$ws.Range("H137").Value = 9739.326999999999
$ws.Range("I137").Value = 4664.081
$ws.Range("J137").Value = 18681.428
$ws.Range("K137").Value = 13992.243
$ws.Range("L137").Value = 56044.284
$ws.Range("M137").Value = -11442.243
$ws.Range("N137").Value = -61144.284

$ws.Range("H138").Value = 4380.2905
$ws.Range("J138").Value = 5245.1304
$ws.Range("L138").Value = 15735.3912
$ws.Range("N138").Value = -26015.3912

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1653.31
$ws.Range("I32").Value = 1584.4183
$ws.Range("J32").Value = 5029
$ws.Range("K32").Value = 1584.4183
$ws.Range("L32").Value = 5029
$ws.Range("M32").Value = -1297.4183
$ws.Range("N32").Value = -5603

$ws.Range("H45").Value = 10672.4
$ws.Range("I45").Value = 11699
$ws.Range("K45").Value = 11699
$ws.Range("M45").Value = -11322

$ws.Range("H61").Value = 7453.1875
$ws.Range("I61").Value = 2954.1738
$ws.Range("K61").Value = 2954.1738
$ws.Range("M61").Value = -2742.1738

$ws.Range("H74").Value = 8163.1523
$ws.Range("I74").Value = 9000.605
$ws.Range("K74").Value = 9000.605
$ws.Range("M74").Value = -8126.605

$ws.Range("H77").Value = 8163.1523
$ws.Range("I77").Value = 9000.605
$ws.Range("K77").Value = 45003.02499999999
$ws.Range("M77").Value = -40635.02499999999

$ws.Range("H88").Value = 2563.9
$ws.Range("I88").Value = 1698.5
$ws.Range("K88").Value = 1698.5
$ws.Range("M88").Value = -1292.5

$ws.Range("H91").Value = 2563.9
$ws.Range("I91").Value = 1698.5
$ws.Range("K91").Value = 1698.5
$ws.Range("M91").Value = -294.5

$ws.Range("H132").Value = 3710.0408
$ws.Range("I132").Value = 2949.3333
$ws.Range("J132").Value = 6676.8
$ws.Range("K132").Value = 8847.999899999999
$ws.Range("L132").Value = 20030.4
$ws.Range("M132").Value = -6317.999899999999
$ws.Range("N132").Value = -25090.4

$ws.Range("H136").Value = 7453.1875
$ws.Range("I136").Value = 2954.1738
$ws.Range("K136").Value = 8862.5214
$ws.Range("M136").Value = -6312.5214

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 15041.963
$ws.Range("I20").Value = 20879.883
$ws.Range("J20").Value = 5117.5
$ws.Range("K20").Value = 20879.883
$ws.Range("L20").Value = 5117.5
$ws.Range("M20").Value = -20632.883
$ws.Range("N20").Value = -5611.5

$ws.Range("H99").Value = 2133.697
$ws.Range("I99").Value = 2106.625
$ws.Range("J99").Value = 3000
$ws.Range("K99").Value = 2106.625
$ws.Range("L99").Value = 3000
$ws.Range("M99").Value = -608.625
$ws.Range("N99").Value = -5996

$ws.Range("H134").Value = 8663.049999999999
$ws.Range("I134").Value = 4927.36
$ws.Range("J134").Value = 14889.2
$ws.Range("K134").Value = 14782.08
$ws.Range("L134").Value = 44667.60000000001
$ws.Range("M134").Value = -12247.08
$ws.Range("N134").Value = -49737.60000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 6256
$ws.Range("J25").Value = 10006.5
$ws.Range("L25").Value = 10006.5
$ws.Range("N25").Value = -10354.5

$ws.Range("H31").Value = 358844.53
$ws.Range("I31").Value = 169394.83
$ws.Range("J31").Value = 404312.44
$ws.Range("K31").Value = 169394.83
$ws.Range("L31").Value = 404312.44
$ws.Range("M31").Value = -169099.83
$ws.Range("N31").Value = -404902.44

$ws.Range("H34").Value = 358844.53
$ws.Range("I34").Value = 169394.83
$ws.Range("J34").Value = 404312.44
$ws.Range("K34").Value = 169394.83
$ws.Range("L34").Value = 404312.44
$ws.Range("M34").Value = -169192.83
$ws.Range("N34").Value = -404716.44

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 2703.875
$ws.Range("I14").Value = 2703.875
$ws.Range("K14").Value = 8111.625
$ws.Range("M14").Value = -7938.625

$ws.Range("H39").Value = 6679.4
$ws.Range("I39").Value = 1425
$ws.Range("J39").Value = 7993
$ws.Range("K39").Value = 4275
$ws.Range("L39").Value = 23979
$ws.Range("M39").Value = -3981
$ws.Range("N39").Value = -24567

$ws.Range("H97").Value = 1773.4445
$ws.Range("I97").Value = 413.77777
$ws.Range("J97").Value = 3133.111
$ws.Range("K97").Value = 1241.33331
$ws.Range("L97").Value = 9399.332999999999
$ws.Range("M97").Value = -745.33331
$ws.Range("N97").Value = -10391.333

$ws.Range("H131").Value = 4666.2383
$ws.Range("J131").Value = 5210.4863
$ws.Range("L131").Value = 15631.4589
$ws.Range("N131").Value = -25711.4589

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H94").Value = 49447
$ws.Range("J94").Value = 49447
$ws.Range("L94").Value = 49447
$ws.Range("N94").Value = -50799

$ws.Range("H132").Value = 11045.433
$ws.Range("I132").Value = 9816.357
$ws.Range("K132").Value = 29449.071
$ws.Range("M132").Value = -26919.071

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H26").Value = 9996.5
$ws.Range("I26").Value = 9995
$ws.Range("K26").Value = 9995
$ws.Range("M26").Value = -9700

$ws.Range("H93").Value = 4474.25
$ws.Range("I93").Value = 4496.857
$ws.Range("J93").Value = 4316
$ws.Range("K93").Value = 4496.857
$ws.Range("L93").Value = 4316
$ws.Range("M93").Value = -3248.857
$ws.Range("N93").Value = -6812

$ws.Range("H136").Value = 2451.373
$ws.Range("I136").Value = 2020.08
$ws.Range("J136").Value = 3719.8823
$ws.Range("K136").Value = 6060.24
$ws.Range("L136").Value = 11159.6469
$ws.Range("M136").Value = -3510.24
$ws.Range("N136").Value = -16259.6469

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 10000
$ws.Range("J18").Value = 10000
$ws.Range("L18").Value = 10000
$ws.Range("N18").Value = -10346

$ws.Range("H100").Value = 714.93335
$ws.Range("I100").Value = 714.93335
$ws.Range("K100").Value = 1429.8667
$ws.Range("M100").Value = -888.8667

$ws.Range("H122").Value = 7227.0713
$ws.Range("I122").Value = 5021
$ws.Range("K122").Value = 15063
$ws.Range("M122").Value = -12613

$ws.Range("H132").Value = 26597.975
$ws.Range("I132").Value = 26035.018
$ws.Range("K132").Value = 78105.054
$ws.Range("M132").Value = -75575.054
